$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row 1 with two new columns (P, Q), copying the style of O1
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Update the data block B2:Q25 (existing columns change values; P & Q are new, filled with 0)
$data = New-Object 'object[,]' 24,16
$data[0,0] = 3.378213651322142
$data[0,1] = 0.5742940435717401
$data[0,2] = 0.03522866393421253
$data[0,3] = 0.03738537503279371
$data[0,4] = 1.121892043285527
$data[0,5] = 0.9002394968517535
$data[0,6] = 0.003202048930136092
$data[0,7] = 0
$data[0,8] = 0.6157497418812881
$data[0,9] = 0.5314621174801104
$data[0,10] = 0.05566840826468633
$data[0,11] = 0
$data[0,12] = 0.3565997824003091
$data[0,13] = 0
$data[0,14] = 0
$data[0,15] = 0
$data[1,0] = 2.942196971035571
$data[1,1] = 0.5018033520778147
$data[1,2] = 0.03189349820311094
$data[1,3] = 0.03510888316151384
$data[1,4] = 1.068080225612803
$data[1,5] = 0.8498243765501314
$data[1,6] = 0.005814421139721915
$data[1,7] = 0
$data[1,8] = 0.5993677995312936
$data[1,9] = 0.5294122098163783
$data[1,10] = 0.0529585807501225
$data[1,11] = 0
$data[1,12] = 0.3143348947597815
$data[1,13] = 0
$data[1,14] = 0
$data[1,15] = 0
$data[2,0] = 2.674967090450423
$data[2,1] = 0.4573441629894717
$data[2,2] = 0.02985488272019765
$data[2,3] = 0.0337099917402135
$data[2,4] = 1.036719017932519
$data[2,5] = 0.8205453378904082
$data[2,6] = 0.007867353992464965
$data[2,7] = 0
$data[2,8] = 0.5902719434921124
$data[2,9] = 0.5291851718132179
$data[2,10] = 0.05127692595167144
$data[2,11] = 0
$data[2,12] = 0.2884326187289616
$data[2,13] = 0
$data[2,14] = 0
$data[2,15] = 0
$data[3,0] = 2.566152532369301
$data[3,1] = 0.4392317726980366
$data[3,2] = 0.02902563783159806
$data[3,3] = 0.03313884530603328
$data[3,4] = 1.024336836035438
$data[3,5] = 0.809008630237841
$data[3,6] = 0.008810733441820665
$data[3,7] = 0
$data[3,8] = 0.5867951527077935
$data[3,9] = 0.5293395876340341
$data[3,10] = 0.0505862425447603
$data[3,11] = 0
$data[3,12] = 0.2778847416932422
$data[3,13] = 0
$data[3,14] = 0
$data[3,15] = 0
$data[4,0] = 2.548087945449026
$data[4,1] = 0.4362243199164197
$data[4,2] = 0.02888801099327054
$data[4,3] = 0.03304391852669042
$data[4,4] = 1.022304113904276
$data[4,5] = 0.8071160546688816
$data[4,6] = 0.008973658302325632
$data[4,7] = 0
$data[4,8] = 0.5862313702911592
$data[4,9] = 0.5293797894903136
$data[4,10] = 0.05047120332032984
$data[4,11] = 0
$data[4,12] = 0.2761335964120804
$data[4,13] = 0
$data[4,14] = 0
$data[4,15] = 0
$data[5,0] = 2.673499288907067
$data[5,1] = 0.4570998819903878
$data[5,2] = 0.02984369415050736
$data[5,3] = 0.03370229449282292
$data[5,4] = 1.036550448503135
$data[5,5] = 0.8203881858503479
$data[5,6] = 0.007879652175801555
$data[5,7] = 0
$data[5,8] = 0.5902241393621068
$data[5,9] = 0.5291862705202028
$data[5,10] = 0.05126763417703284
$data[5,11] = 0
$data[5,12] = 0.2882903419651655
$data[5,13] = 0
$data[5,14] = 0
$data[5,15] = 0
$data[6,0] = 3.227750101119511
$data[6,1] = 0.5492840475072853
$data[6,2] = 0.03407628899776682
$data[6,3] = 0.03660017225611956
$data[6,4] = 1.102973971781466
$data[6,5] = 0.8824919729382543
$data[6,6] = 0.004005860460893018
$data[6,7] = 0
$data[6,8] = 0.6098938645172609
$data[6,9] = 0.5305334967442548
$data[6,10] = 0.05473718854871112
$data[6,11] = 0
$data[6,12] = 0.3420138524418519
$data[6,13] = 0
$data[6,14] = 0
$data[6,15] = 0
$data[7,0] = 4.320415876822665
$data[7,1] = 0.730834515632381
$data[7,2] = 0.04248895769599059
$data[7,3] = 0.04231298132938122
$data[7,4] = 1.247734538434898
$data[7,5] = 1.018858210989023
$data[7,6] = 0.0002612579403886706
$data[7,7] = 0
$data[7,8] = 0.656692037972519
$data[7,9] = 0.5419557498354877
$data[7,10] = 0.06144371126460779
$data[7,11] = 0
$data[7,12] = 0.4479812328970638
$data[7,13] = 0
$data[7,14] = 0
$data[7,15] = 0
$data[8,0] = 5.128952017444362
$data[8,1] = 0.869059233120197
$data[8,2] = 0.04834724044744121
$data[8,3] = 0.04559592303243676
$data[8,4] = 1.340469153859829
$data[8,5] = 1.107856760170037
$data[8,6] = 0.0004302313997999185
$data[8,7] = 0
$data[8,8] = 0.6859164802244493
$data[8,9] = 0.5446424802941365
$data[8,10] = 0.06547191315135015
$data[8,11] = 0
$data[8,12] = 0.5148382078785829
$data[8,13] = 0
$data[8,14] = 0
$data[8,15] = 0
$data[9,0] = 5.492313975234765
$data[9,1] = 0.9665920742392302
$data[9,2] = 0.04717089501021121
$data[9,3] = 0.04070858670855948
$data[9,4] = 1.181298931409344
$data[9,5] = 0.9662501012007567
$data[9,6] = 0.01915232211114315
$data[9,7] = 0
$data[9,8] = 0.6076299435575123
$data[9,9] = 0.4497569450123393
$data[9,10] = 0.06558996054978561
$data[9,11] = 0
$data[9,12] = 0.4472294367971017
$data[9,13] = 0
$data[9,14] = 0
$data[9,15] = 0
$data[10,0] = 5.627550453819538
$data[10,1] = 1.021557117770698
$data[10,2] = 0.04497126317965439
$data[10,3] = 0.03886028513967543
$data[10,4] = 1.034563680391784
$data[10,5] = 0.8348909168231984
$data[10,6] = 0.05779340307736192
$data[10,7] = 0
$data[10,8] = 0.5385278100825985
$data[10,9] = 0.3763870619622622
$data[10,10] = 0.07318351214461671
$data[10,11] = 0
$data[10,12] = 0.3807950761386252
$data[10,13] = 0
$data[10,14] = 0
$data[10,15] = 0
$data[11,0] = 5.593468306582167
$data[11,1] = 1.046265524444522
$data[11,2] = 0.04185753007564585
$data[11,3] = 0.03918202244606839
$data[11,4] = 0.8888420570317379
$data[11,5] = 0.7036064752890212
$data[11,6] = 0.1133533095533608
$data[11,7] = 0
$data[11,8] = 0.4723005615089733
$data[11,9] = 0.3148947670537972
$data[11,10] = 0.0867036501337104
$data[11,11] = 0
$data[11,12] = 0.3119831365791583
$data[11,13] = 0
$data[11,14] = 0
$data[11,15] = 0
$data[12,0] = 5.493064648802942
$data[12,1] = 1.049845029037499
$data[12,2] = 0.03929323680396024
$data[12,3] = 0.04075634159942876
$data[12,4] = 0.7892570135586183
$data[12,5] = 0.6135748465985955
$data[12,6] = 0.1625458399557829
$data[12,7] = 0
$data[12,8] = 0.4280862550198918
$data[12,9] = 0.2784779919350591
$data[12,10] = 0.09977247945621315
$data[12,11] = 0
$data[12,12] = 0.2639373681295893
$data[12,13] = 0
$data[12,14] = 0
$data[12,15] = 0
$data[13,0] = 5.432566806173895
$data[13,1] = 1.043935545164118
$data[13,2] = 0.0384761844368775
$data[13,3] = 0.0412040163389058
$data[13,4] = 0.7644503223613555
$data[13,5] = 0.5908996330676786
$data[13,6] = 0.1749976013633727
$data[13,7] = 0
$data[13,8] = 0.4175807189980674
$data[13,9] = 0.2712060164902113
$data[13,10] = 0.1030702055947259
$data[13,11] = 0
$data[13,12] = 0.2512541370982291
$data[13,13] = 0
$data[13,14] = 0
$data[13,15] = 0
$data[14,0] = 5.089608835301931
$data[14,1] = 0.9778304404864002
$data[14,2] = 0.03664324262283714
$data[14,3] = 0.0396626203647612
$data[14,4] = 0.7553969965176677
$data[14,5] = 0.5800412293465342
$data[14,6] = 0.1621161439121863
$data[14,7] = 0
$data[14,8] = 0.4185467993785181
$data[14,9] = 0.282014913603188
$data[14,10] = 0.09803130491518175
$data[14,11] = 0
$data[14,12] = 0.2388283322743945
$data[14,13] = 0
$data[14,14] = 0
$data[14,15] = 0
$data[15,0] = 4.880764911542713
$data[15,1] = 0.9254795391275366
$data[15,2] = 0.03658643875036915
$data[15,3] = 0.03740220609595735
$data[15,4] = 0.8005975517199104
$data[15,5] = 0.6188954257335553
$data[15,6] = 0.1243782984726209
$data[15,7] = 0
$data[15,8] = 0.4425445910693071
$data[15,9] = 0.3092066305222687
$data[15,10] = 0.08662815822889058
$data[15,11] = 0
$data[15,12] = 0.2539376786437515
$data[15,13] = 0
$data[15,14] = 0
$data[15,15] = 0
$data[16,0] = 4.762408207456645
$data[16,1] = 0.8795471382317146
$data[16,2] = 0.03807381381362163
$data[16,3] = 0.0356909433239303
$data[16,4] = 0.9016032510814398
$data[16,5] = 0.7086972272039702
$data[16,6] = 0.0715919454901055
$data[16,7] = 0
$data[16,8] = 0.4908121313088145
$data[16,9] = 0.3572702759766422
$data[16,10] = 0.07238416923015478
$data[16,11] = 0
$data[16,12] = 0.2971263022092643
$data[16,13] = 0
$data[16,14] = 0
$data[16,15] = 0
$data[17,0] = 4.724984446257452
$data[17,1] = 0.842918594909861
$data[17,2] = 0.04076558817153852
$data[17,3] = 0.03666069672710615
$data[17,4] = 1.047491265620877
$data[17,5] = 0.8399195756625772
$data[17,6] = 0.0263071488080584
$data[17,7] = 0
$data[17,8] = 0.5583238097097194
$data[17,9] = 0.4246757511749877
$data[17,10] = 0.06262348126841744
$data[17,11] = 0
$data[17,12] = 0.3642568982104706
$data[17,13] = 0
$data[17,14] = 0
$data[17,15] = 0
$data[18,0] = 4.916074992141773
$data[18,1] = 0.8328153481903371
$data[18,2] = 0.04677898940100533
$data[18,3] = 0.04468232560508767
$data[18,4] = 1.314328764814221
$data[18,5] = 1.08276499272749
$data[18,6] = 0.0001710473031542037
$data[18,7] = 0
$data[18,8] = 0.6773456261008874
$data[18,9] = 0.5429674981614312
$data[18,10] = 0.06436342883439616
$data[18,11] = 0
$data[18,12] = 0.4967532109524484
$data[18,13] = 0
$data[18,14] = 0
$data[18,15] = 0
$data[19,0] = 5.540536527318864
$data[19,1] = 0.9335470363071465
$data[19,2] = 0.05203290446307562
$data[19,3] = 0.04877916035693719
$data[19,4] = 1.426819411599567
$data[19,5] = 1.189060198388745
$data[19,6] = 0.001105925369503824
$data[19,7] = 0
$data[19,8] = 0.7189888065259993
$data[19,9] = 0.5657514951922451
$data[19,10] = 0.06888705541336027
$data[19,11] = 0
$data[19,12] = 0.566492844848355
$data[19,13] = 0
$data[19,14] = 0
$data[19,15] = 0
$data[20,0] = 5.950508840100156
$data[20,1] = 1.001696876532208
$data[20,2] = 0.05528207245099281
$data[20,3] = 0.05098487333618973
$data[20,4] = 1.490544436893728
$data[20,5] = 1.249960298709425
$data[20,6] = 0.002433069706141877
$data[20,7] = 0
$data[20,8] = 0.7419520097922145
$data[20,9] = 0.5759918634380696
$data[20,10] = 0.07140012082281455
$data[20,11] = 0
$data[20,12] = 0.6063692302128203
$data[20,13] = 0
$data[20,14] = 0
$data[20,15] = 0
$data[21,0] = 5.731474873491209
$data[21,1] = 0.9652835073923711
$data[21,2] = 0.05354332536192885
$data[21,3] = 0.04980407078817972
$data[21,4] = 1.456289397033586
$data[21,5] = 1.217203439218167
$data[21,6] = 0.001662064369820415
$data[21,7] = 0
$data[21,8] = 0.7295643724483796
$data[21,9] = 0.5703887479418199
$data[21,10] = 0.07005618666923752
$data[21,11] = 0
$data[21,12] = 0.5850609529945956
$data[21,13] = 0
$data[21,14] = 0
$data[21,15] = 0
$data[22,0] = 4.906452183067756
$data[22,1] = 0.8281861302947107
$data[22,2] = 0.04705074057669378
$data[22,3] = 0.04540253437216801
$data[22,4] = 1.331680213904164
$data[22,5] = 1.098456109662038
$data[22,6] = 0.00006253031109326734
$data[22,7] = 0
$data[22,8] = 0.6854270768593267
$data[22,9] = 0.552075537000313
$data[22,10] = 0.06501599627570442
$data[22,11] = 0
$data[22,12] = 0.5048746285184649
$data[22,13] = 0
$data[22,14] = 0
$data[22,15] = 0
$data[23,0] = 4.023909773170885
$data[23,1] = 0.681579929562929
$data[23,2] = 0.0401951239580427
$data[23,3] = 0.04075854060905293
$data[23,4] = 1.206846491349225
$data[23,5] = 0.980217655339743
$data[23,6] = 0.0008195009018061583
$data[23,7] = 0
$data[23,8] = 0.6430648036778734
$data[23,9] = 0.5378384983711371
$data[23,10] = 0.0596333419773547
$data[23,11] = 0
$data[23,12] = 0.4192137219651926
$data[23,13] = 0
$data[23,14] = 0
$data[23,15] = 0

$ws.Range("B2:Q25").Value = $data
